# Add a new "register" test-case row (row 2) below the existing login row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "firstname"
$ws.Range("B2").Value = "first@gmail"

# Turn the email in B2 into a mailto hyperlink, styled like the existing
# hyperlink in A1 (re-applies the built-in "Hyperlink" cell style).
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:first@gmail")
$ws.Range("B2").Style = "Hyperlink"

# Matches the saved selection state (active cell B2) in the target file.
$ws.Range("B2").Select()
